$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.161539316177368
$ws.Range("B1").Value = 2.415287017822266
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.379549503326416
$ws.Range("E1").Value = 1.230799913406372
